# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 93
$ws.Range("A93").Value = 91
$ws.Range("B93").Value = 6236611
$ws.Range("C93").Value = "Venezuela Primera Division"
$ws.Range("D93").Value = "Venezuela Primera Division"
$ws.Range("E93").Value = 45199.6875
$ws.Range("F93").Value = "Mineros"
$ws.Range("G93").Value = "Monagas"
$ws.Range("H93").Value = 2
$ws.Range("I93").Value = 1
$ws.Range("J93").Value = "H"
$ws.Range("K93").Value = 3.2
$ws.Range("L93").Value = 3.4
$ws.Range("M93").Value = 2
$ws.Range("N93").Value = 4.2
$ws.Range("O93").Value = 3.8
$ws.Range("P93").Value = 1.65
$ws.Range("Q93").Value = 0.75
$ws.Range("R93").Value = 1.95
$ws.Range("S93").Value = 1.85
$ws.Range("T93").Value = 2.5
$ws.Range("U93").Value = 1.825
$ws.Range("V93").Value = 1.975
$ws.Range("W93").Value = 3.2
$ws.Range("X93").Value = -1
$ws.Range("Y93").Value = -1
$ws.Range("Z93").Value = 0.95
$ws.Range("AA93").Value = -1
$ws.Range("AB93").Value = 0.825
$ws.Range("AC93").Value = -1

# Row 94
$ws.Range("A94").Value = 92
$ws.Range("B94").Value = 6236254
$ws.Range("C94").Value = "Venezuela Primera Division"
$ws.Range("D94").Value = "Venezuela Primera Division"
$ws.Range("E94").Value = 45199.6875
$ws.Range("F94").Value = "Academia Puerto Cabello"
$ws.Range("G94").Value = "Estudiantes Merida"
$ws.Range("H94").Value = 1
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = "H"
$ws.Range("K94").Value = 1.727
$ws.Range("L94").Value = 3.4
$ws.Range("M94").Value = 4.333
$ws.Range("N94").Value = 1.666
$ws.Range("O94").Value = 3.4
$ws.Range("P94").Value = 4.75
$ws.Range("Q94").Value = -0.75
$ws.Range("R94").Value = 1.875
$ws.Range("S94").Value = 1.925
$ws.Range("T94").Value = 2.5
$ws.Range("U94").Value = 1.9
$ws.Range("V94").Value = 1.9
$ws.Range("W94").Value = 0.6659999999999999
$ws.Range("X94").Value = -1
$ws.Range("Y94").Value = -1
$ws.Range("Z94").Value = 0.4375
$ws.Range("AA94").Value = -0.5
$ws.Range("AB94").Value = -1
$ws.Range("AC94").Value = 0.8999999999999999

# Row 95
$ws.Range("A95").Value = 93
$ws.Range("B95").Value = 6236253
$ws.Range("C95").Value = "Venezuela Primera Division"
$ws.Range("D95").Value = "Venezuela Primera Division"
$ws.Range("E95").Value = 45199.6875
$ws.Range("F95").Value = "Deportivo La Guaira"
$ws.Range("G95").Value = "UCV"
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = "D"
$ws.Range("K95").Value = 1.833
$ws.Range("L95").Value = 3.25
$ws.Range("M95").Value = 4
$ws.Range("N95").Value = 2
$ws.Range("O95").Value = 3.2
$ws.Range("P95").Value = 3.5
$ws.Range("Q95").Value = -0.25
$ws.Range("R95").Value = 1.775
$ws.Range("S95").Value = 2.025
$ws.Range("T95").Value = 2.25
$ws.Range("U95").Value = 1.9
$ws.Range("V95").Value = 1.9
$ws.Range("W95").Value = -1
$ws.Range("X95").Value = 2.2
$ws.Range("Y95").Value = -1
$ws.Range("Z95").Value = -0.5
$ws.Range("AA95").Value = 0.5125
$ws.Range("AB95").Value = -1
$ws.Range("AC95").Value = 0.8999999999999999

# Row 96
$ws.Range("A96").Value = 94
$ws.Range("B96").Value = 6236252
$ws.Range("C96").Value = "Venezuela Primera Division"
$ws.Range("D96").Value = "Venezuela Primera Division"
$ws.Range("E96").Value = 45199.6875
$ws.Range("F96").Value = "Deportivo Tachira"
$ws.Range("G96").Value = "CD Hermanos Colmenares"
$ws.Range("H96").Value = 1
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = "H"
$ws.Range("K96").Value = 1.363
$ws.Range("L96").Value = 4.2
$ws.Range("M96").Value = 7.5
$ws.Range("N96").Value = 1.333
$ws.Range("O96").Value = 4.5
$ws.Range("P96").Value = 8
$ws.Range("Q96").Value = -1.5
$ws.Range("R96").Value = 2
$ws.Range("S96").Value = 1.8
$ws.Range("T96").Value = 2.5
$ws.Range("U96").Value = 1.925
$ws.Range("V96").Value = 1.875
$ws.Range("W96").Value = 0.333
$ws.Range("X96").Value = -1
$ws.Range("Y96").Value = -1
$ws.Range("Z96").Value = -1
$ws.Range("AA96").Value = 0.8
$ws.Range("AB96").Value = -1
$ws.Range("AC96").Value = 0.875

# Row 97
$ws.Range("A97").Value = 95
$ws.Range("B97").Value = 6236255
$ws.Range("C97").Value = "Venezuela Primera Division"
$ws.Range("D97").Value = "Venezuela Primera Division"
$ws.Range("E97").Value = 45199.6875
$ws.Range("F97").Value = "Deportivo Rayo Zuliano"
$ws.Range("G97").Value = "Caracas"
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = "D"
$ws.Range("K97").Value = 3.75
$ws.Range("L97").Value = 3.1
$ws.Range("M97").Value = 1.95
$ws.Range("N97").Value = 2.9
$ws.Range("O97").Value = 2.875
$ws.Range("P97").Value = 2.45
$ws.Range("Q97").Value = 0.25
$ws.Range("R97").Value = 1.775
$ws.Range("S97").Value = 2.025
$ws.Range("T97").Value = 2.25
$ws.Range("U97").Value = 1.85
$ws.Range("V97").Value = 1.95
$ws.Range("W97").Value = -1
$ws.Range("X97").Value = 1.875
$ws.Range("Y97").Value = -1
$ws.Range("Z97").Value = 0.3875
$ws.Range("AA97").Value = -0.5
$ws.Range("AB97").Value = -1
$ws.Range("AC97").Value = 0.95

# Row 98
$ws.Range("A98").Value = 96
$ws.Range("B98").Value = 6236612
$ws.Range("C98").Value = "Venezuela Primera Division"
$ws.Range("D98").Value = "Venezuela Primera Division"
$ws.Range("E98").Value = 45199.6875
$ws.Range("F98").Value = "Zamora"
$ws.Range("G98").Value = "Carabobo"
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 2
$ws.Range("J98").Value = "A"
$ws.Range("K98").Value = 3.2
$ws.Range("L98").Value = 3.1
$ws.Range("M98").Value = 2.15
$ws.Range("N98").Value = 4.5
$ws.Range("O98").Value = 3.3
$ws.Range("P98").Value = 1.75
$ws.Range("Q98").Value = 0.5
$ws.Range("R98").Value = 2
$ws.Range("S98").Value = 1.8
$ws.Range("T98").Value = 2.25
$ws.Range("U98").Value = 1.925
$ws.Range("V98").Value = 1.875
$ws.Range("W98").Value = -1
$ws.Range("X98").Value = -1
$ws.Range("Y98").Value = 0.75
$ws.Range("Z98").Value = -1
$ws.Range("AA98").Value = 0.8
$ws.Range("AB98").Value = -0.5
$ws.Range("AC98").Value = 0.4375

# Row 135
$ws.Range("A135").Value = 133
$ws.Range("B135").Value = 7842507
$ws.Range("C135").Value = "Venezuela Primera Division"
$ws.Range("D135").Value = "Venezuela Primera Division"
$ws.Range("E135").Value = 45339.78125
$ws.Range("F135").Value = "Academia Puerto Cabello"
$ws.Range("G135").Value = "Estudiantes Merida"
$ws.Range("H135").Value = 2
$ws.Range("I135").Value = 1
$ws.Range("J135").Value = "H"
$ws.Range("K135").Value = 1.727
$ws.Range("L135").Value = 3.5
$ws.Range("M135").Value = 4.2
$ws.Range("N135").Value = 1.85
$ws.Range("O135").Value = 3.5
$ws.Range("P135").Value = 3.6
$ws.Range("Q135").Value = -0.5
$ws.Range("R135").Value = 1.925
$ws.Range("S135").Value = 1.875
$ws.Range("T135").Value = 2.5
$ws.Range("U135").Value = 1.9
$ws.Range("V135").Value = 1.9
$ws.Range("W135").Value = 0.8500000000000001
$ws.Range("X135").Value = -1
$ws.Range("Y135").Value = -1
$ws.Range("Z135").Value = 0.925
$ws.Range("AA135").Value = -1
$ws.Range("AB135").Value = 0.8999999999999999
$ws.Range("AC135").Value = -1

# Row 136
$ws.Range("A136").Value = 134
$ws.Range("B136").Value = 7842504
$ws.Range("C136").Value = "Venezuela Primera Division"
$ws.Range("D136").Value = "Venezuela Primera Division"
$ws.Range("E136").Value = 45339.78125
$ws.Range("F136").Value = "Angostura FC"
$ws.Range("G136").Value = "Deportivo La Guaira"
$ws.Range("H136").Value = 1
$ws.Range("I136").Value = 1
$ws.Range("J136").Value = "D"
$ws.Range("K136").Value = 2.75
$ws.Range("L136").Value = 3
$ws.Range("M136").Value = 2.45
$ws.Range("N136").Value = 3.1
$ws.Range("O136").Value = 2.875
$ws.Range("P136").Value = 2.3
$ws.Range("Q136").Value = 0.25
$ws.Range("R136").Value = 1.8
$ws.Range("S136").Value = 2
$ws.Range("T136").Value = 2.25
$ws.Range("U136").Value = 2.05
$ws.Range("V136").Value = 1.75
$ws.Range("W136").Value = -1
$ws.Range("X136").Value = 1.875
$ws.Range("Y136").Value = -1
$ws.Range("Z136").Value = 0.4
$ws.Range("AA136").Value = -0.5
$ws.Range("AB136").Value = -0.5
$ws.Range("AC136").Value = 0.375

# Row 157
$ws.Range("A157").Value = 155
$ws.Range("B157").Value = 7920998
$ws.Range("C157").Value = "Venezuela Primera Division"
$ws.Range("D157").Value = "Venezuela Primera Division"
$ws.Range("E157").Value = 45360.79166666666
$ws.Range("F157").Value = "Zamora"
$ws.Range("G157").Value = "Caracas"
$ws.Range("H157").Value = 2
$ws.Range("I157").Value = 2
$ws.Range("J157").Value = "D"
$ws.Range("K157").Value = 3.75
$ws.Range("L157").Value = 3.2
$ws.Range("M157").Value = 1.909
$ws.Range("N157").Value = 3
$ws.Range("O157").Value = 2.9
$ws.Range("P157").Value = 2.375
$ws.Range("Q157").Value = 0.25
$ws.Range("R157").Value = 1.8
$ws.Range("S157").Value = 2
$ws.Range("T157").Value = 2
$ws.Range("U157").Value = 1.825
$ws.Range("V157").Value = 1.975
$ws.Range("W157").Value = -1
$ws.Range("X157").Value = 1.9
$ws.Range("Y157").Value = -1
$ws.Range("Z157").Value = 0.4
$ws.Range("AA157").Value = -0.5
$ws.Range("AB157").Value = 0.825
$ws.Range("AC157").Value = -1

# Row 158
$ws.Range("A158").Value = 156
$ws.Range("B158").Value = 7920997
$ws.Range("C158").Value = "Venezuela Primera Division"
$ws.Range("D158").Value = "Venezuela Primera Division"
$ws.Range("E158").Value = 45360.79166666666
$ws.Range("F158").Value = "Carabobo"
$ws.Range("G158").Value = "UCV"
$ws.Range("H158").Value = 0
$ws.Range("I158").Value = 1
$ws.Range("J158").Value = "A"
$ws.Range("K158").Value = 1.833
$ws.Range("L158").Value = 3.1
$ws.Range("M158").Value = 4.2
$ws.Range("N158").Value = 1.833
$ws.Range("O158").Value = 3.1
$ws.Range("P158").Value = 4.2
$ws.Range("Q158").Value = -0.5
$ws.Range("R158").Value = 1.9
$ws.Range("S158").Value = 1.9
$ws.Range("T158").Value = 2
$ws.Range("U158").Value = 1.85
$ws.Range("V158").Value = 1.95
$ws.Range("W158").Value = -1
$ws.Range("X158").Value = -1
$ws.Range("Y158").Value = 3.2
$ws.Range("Z158").Value = -1
$ws.Range("AA158").Value = 0.8999999999999999
$ws.Range("AB158").Value = -1
$ws.Range("AC158").Value = 0.95

# Row 173
$ws.Range("A173").Value = 171
$ws.Range("B173").Value = 7958192
$ws.Range("C173").Value = "Venezuela Primera Division"
$ws.Range("D173").Value = "Venezuela Primera Division"
$ws.Range("E173").Value = 45371.89583333334
$ws.Range("F173").Value = "Deportivo Tachira"
$ws.Range("G173").Value = "Monagas"
$ws.Range("H173").Value = 1
$ws.Range("I173").Value = 0
$ws.Range("J173").Value = "H"
$ws.Range("K173").Value = 1.666
$ws.Range("L173").Value = 3.4
$ws.Range("M173").Value = 4.5
$ws.Range("N173").Value = 1.95
$ws.Range("O173").Value = 3.25
$ws.Range("P173").Value = 3.5
$ws.Range("Q173").Value = -0.5
$ws.Range("R173").Value = 1.975
$ws.Range("S173").Value = 1.825
$ws.Range("T173").Value = 2.25
$ws.Range("U173").Value = 2.025
$ws.Range("V173").Value = 1.775
$ws.Range("W173").Value = 0.95
$ws.Range("X173").Value = -1
$ws.Range("Y173").Value = -1
$ws.Range("Z173").Value = 0.9750000000000001
$ws.Range("AA173").Value = -1
$ws.Range("AB173").Value = -1
$ws.Range("AC173").Value = 0.7749999999999999

# Row 174
$ws.Range("A174").Value = 172
$ws.Range("B174").Value = 7958193
$ws.Range("C174").Value = "Venezuela Primera Division"
$ws.Range("D174").Value = "Venezuela Primera Division"
$ws.Range("E174").Value = 45371.89583333334
$ws.Range("F174").Value = "Zamora"
$ws.Range("G174").Value = "Academia Puerto Cabello"
$ws.Range("H174").Value = 0
$ws.Range("I174").Value = 0
$ws.Range("J174").Value = "D"
$ws.Range("K174").Value = 3.75
$ws.Range("L174").Value = 3.3
$ws.Range("M174").Value = 1.85
$ws.Range("N174").Value = 3.1
$ws.Range("O174").Value = 3.2
$ws.Range("P174").Value = 2.1
$ws.Range("Q174").Value = 0.25
$ws.Range("R174").Value = 1.875
$ws.Range("S174").Value = 1.925
$ws.Range("T174").Value = 2.25
$ws.Range("U174").Value = 2.025
$ws.Range("V174").Value = 1.775
$ws.Range("W174").Value = -1
$ws.Range("X174").Value = 2.2
$ws.Range("Y174").Value = -1
$ws.Range("Z174").Value = 0.4375
$ws.Range("AA174").Value = -0.5
$ws.Range("AB174").Value = -1
$ws.Range("AC174").Value = 0.7749999999999999

# New row 196
$ws.Range("A196").Value = 194
$ws.Range("B196").Value = 7977874
$ws.Range("C196").Value = "Venezuela Primera Division"
$ws.Range("D196").Value = "Venezuela Primera Division"
$ws.Range("E196").Value = 45394.83333333334
$ws.Range("F196").Value = "Monagas"
$ws.Range("G196").Value = "Zamora"
$ws.Range("K196").Value = 1.615
$ws.Range("L196").Value = 3.5
$ws.Range("M196").Value = 4.75
$ws.Range("N196").Value = 1.666
$ws.Range("O196").Value = 3.5
$ws.Range("P196").Value = 4.333
$ws.Range("Q196").Value = -0.75
$ws.Range("R196").Value = 1.875
$ws.Range("S196").Value = 1.925
$ws.Range("T196").Value = 2.5
$ws.Range("U196").Value = 1.8
$ws.Range("V196").Value = 2
$ws.Range("W196").Value = 0
$ws.Range("X196").Value = 0
$ws.Range("Y196").Value = 0
$ws.Range("Z196").Value = 0
$ws.Range("AA196").Value = 0

# Copy styles for new row 196 id/date cells to match existing data rows
$ws.Range("A195").Copy() | Out-Null
$ws.Range("A196").PasteSpecial(-4122) | Out-Null
$ws.Range("E195").Copy() | Out-Null
$ws.Range("E196").PasteSpecial(-4122) | Out-Null